# Apply "Update results with new concentrations" changes across sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "Info": update objective/time result ---
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Range("A2").Value = 640108574274.0112
$wsInfo.Range("B2").Value = 0.6050000190734863

# --- Sheet "Activados": proceso column -> 1, extend data to 19 rows (rows 2-20) ---
$wsAct = $wb.Worksheets.Item("Activados")
$wsAct.Range("A2:A20").Value = 1
for ($i = 0; $i -lt 19; $i++) {
    $row = 2 + $i
    $wsAct.Cells.Item($row, 2).Value = $i * 20
}

# --- Sheet "Operando": proceso column -> 1 for all data rows (Tiempo column unchanged) ---
$wsOp = $wb.Worksheets.Item("Operando")
$wsOp.Range("A2:A366").Value = 1

# --- Sheet "Contaminantes": update mass (B) and concentration (C) columns ---
$wsCont = $wb.Worksheets.Item("Contaminantes")
$wsCont.Range("B2").Value = 449208244800.0004
$wsCont.Range("C2").Value = 16.66000000000001

$wsCont.Range("B3").Value = 13481640000.00001
$wsCont.Range("C3").Value = 0.5000000000000004

$wsCont.Range("B4").Value = 87091394399.99998
$wsCont.Range("C4").Value = 3.23

$wsCont.Range("B5").Value = 307074.010608
$wsCont.Range("C5").Value = 0.0000113886

$wsCont.Range("B6").Value = 90326988000.00008
$wsCont.Range("C6").Value = 3.350000000000003
